$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 11 (shifts existing rows 11-38 down to 13-40)
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# New row 11 data
$ws.Cells.Item(11, 1).Value = 5
$ws.Cells.Item(11, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(11, 3).Value = "Maule"
$ws.Cells.Item(11, 4).Value = 44536
$ws.Cells.Item(11, 5).Value = 7
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100101
$ws.Cells.Item(11, 8).Value = "Berries"
$ws.Cells.Item(11, 9).Value = 100101001
$ws.Cells.Item(11, 10).Value = "Arándano (blue)"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 180
$ws.Cells.Item(11, 14).Value = 3600
$ws.Cells.Item(11, 15).Value = 3600
$ws.Cells.Item(11, 16).Value = 3600
$ws.Cells.Item(11, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(11, 18).Value = "Provincia de Linares"
$ws.Cells.Item(11, 19).Value = 1800
$ws.Cells.Item(11, 20).Value = 2

# New row 12 data
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Cells.Item(12, 4).Value = 44536
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100101
$ws.Cells.Item(12, 8).Value = "Berries"
$ws.Cells.Item(12, 9).Value = 100101001
$ws.Cells.Item(12, 10).Value = "Arándano (blue)"
$ws.Cells.Item(12, 11).Value = "Sin especificar"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 80
$ws.Cells.Item(12, 14).Value = 3000
$ws.Cells.Item(12, 15).Value = 3000
$ws.Cells.Item(12, 16).Value = 3000
$ws.Cells.Item(12, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(12, 18).Value = "Provincia de Linares"
$ws.Cells.Item(12, 19).Value = 1500
$ws.Cells.Item(12, 20).Value = 2
